# Tasks for Code Reviews.xlsx - update for Code Review 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 ("Bug reporting system set up") - percentage set to 50%
$ws.Range("D33").Value = 0.5

# Row 34 ("Deployment working") - percentage set to 50%
$ws.Range("D34").Value = 0.5

# Row 35 ("PopRep") - not a task of code review 4, so assigned 0%,
# with a note explaining why added to the task description.
$ws.Range("D35").Value = 0
$ws.Range("C35").Value = '"PopRep" (This wasn''t a task of code review 4 so I assigned it 0% for now)'

# Leave the cursor on F31, matching where review work left off.
$ws.Range("F31").Select()
